$wb = $excel.ActiveWorkbook

# Week 5 is the sheet being filled in with new hour entries
$ws = $wb.Worksheets.Item("Week 5")

$ws.Range("A3").Value = "meetings"
$ws.Range("B3").Value = 4

$ws.Range("A4").Value = "presentation"
$ws.Range("B4").Value = 3

$ws.Range("A5").Value = "presentation prep"
$ws.Range("B5").Value = 1

$ws.Range("A6").Value = "implementation"
$ws.Range("B6").Value = 4

# Select cell C6 on this sheet, and make it the active sheet/tab
$ws.Range("C6").Select()
$ws.Activate()
